$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 10
$ws.Range("C2").Value = 40
$ws.Range("B3").Value = 30
$ws.Range("E1").Formula = "=AVERAGE((A1+C2)/2)"

[void]$ws.Range("G2").Select()
